$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $value)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$changes = @(
    @('D2', '43.450.87', $false),
    @('E2', '  +2.55%  ', $false),
    @('D3', '2.313.82', $false),
    @('E3', '  +1.68%  ', $false),
    @('E4', '  -0.04%  ', $false),
    @('D5', '311.11', $true),
    @('E5', '  +0.75%  ', $false),
    @('D6', '104.38', $true),
    @('E6', '  +6.77%  ', $false),
    @('E7', '  +1.14%  ', $false),
    @('E8', '  -0.01%  ', $false),
    @('E9', '  +8.32%  ', $false),
    @('D10', '36.78', $true),
    @('E10', '  +4.52%  ', $false),
    @('D11', '52.81', $true),
    @('E11', '  +1.35%  ', $false),
    @('D12', '0.0815', $true),
    @('E12', '  +0.36%  ', $false),
    @('D13', '0.113', $true),
    @('E13', '  -1.07%  ', $false),
    @('D14', '7.03', $true),
    @('E14', '  +2.66%  ', $false),
    @('D15', '2.671.89', $false),
    @('E15', '  +1.68%  ', $false),
    @('D16', '15.16', $true),
    @('E16', '  +3.37%  ', $false),
    @('D17', '2.318.39', $false),
    @('E17', '  +2.44%  ', $false),
    @('E18', '  +2.48%  ', $false),
    @('D19', '43.358.44', $false),
    @('E19', '  +2.64%  ', $false),
    @('D20', '12.21', $true),
    @('E20', '  -0.75%  ', $false),
    @('D21', '0.0₃0928', $false),
    @('E21', '  +2.30%  ', $false),
    @('E22', '  +3.37%  ', $false),
    @('D23', '68.18', $true),
    @('E23', '  +0.72%  ', $false),
    @('D24', '242.80', $true),
    @('E24', '  +2.51%  ', $false),
    @('E25', '  +2.65%  ', $false),
    @('E26', '  +0.87%  ', $false),
    @('E27', '  +0.06%  ', $false),
    @('D28', '24.90', $true),
    @('E28', '  +5.39%  ', $false),
    @('E29', '  +10.26%  ', $false),
    @('D30', '37.07', $true),
    @('E30', '  -1.53%  ', $false),
    @('E31', '  +0.56%  ', $false),
    @('D32', '168.06', $true),
    @('E32', '  +2.57%  ', $false),
    @('D33', '5.30', $true),
    @('E33', '  +0.80%  ', $false),
    @('E34', '  -0.07%  ', $false),
    @('D35', '18.43', $true),
    @('E35', '  +3.68%  ', $false),
    @('E36', '  +6.72%  ', $false),
    @('D37', '0.0744', $true),
    @('E37', '  +1.13%  ', $false),
    @('D38', '3.06', $true),
    @('E38', '  -1.14%  ', $false),
    @('E39', '  +3.23%  ', $false),
    @('B40', 'Kaspa', $false),
    @('C40', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', $false),
    @('D40', '0.106', $true),
    @('E40', '  +1.99%  ', $false),
    @('B41', 'RenderToken', $false),
    @('C41', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', $false),
    @('D41', '4.48', $true),
    @('E41', '  +7.45%  ', $false),
    @('E42', '  +0.65%  ', $false),
    @('D43', '2.72', $true),
    @('E43', '  +20.00%  ', $false),
    @('E44', '  +3.48%  ', $false),
    @('D45', '1.992.14', $false),
    @('E45', '  +2.26%  ', $false),
    @('D46', '19.05', $true),
    @('D47', '3.06', $true),
    @('E47', '  +2.98%  ', $false),
    @('D48', '10.02', $true),
    @('E48', '  +2.13%  ', $false),
    @('D49', '56.18', $true),
    @('E49', '  +3.28%  ', $false),
    @('D50', '2.95', $true),
    @('E50', '  +2.23%  ', $false),
    @('D51', '1.60', $true),
    @('E51', '  +8.90%  ', $false),
)

foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    $forceText = $chg[2]
    if ($forceText) {
        Set-TextCell $ref $val
    } else {
        $ws.Range($ref).Value = $val
    }
}
